# Update to the "Avverkningsanmälningar" sheet: refresh the "Förändrad" date
# column, reorder a few rows, and append the new row for "A 36808-2023"
# while replacing the old, resolved "A 47107-2023" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the bottom of the table (row 25) so the table grows by one.
$ws.Rows.Item(25).EntireRow.Insert()

# All "Förändrad" (last-checked) values in column C move from 45233 to 45234.
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 3).Value = 45234
}

# Row 20 and 21 swap their "Beteckning" / "Area (ha)" values.
$ws.Cells.Item(20, 1).Value = "A 27505-2023"
$ws.Cells.Item(20, 7).Value = 2.5

$ws.Cells.Item(21, 1).Value = "A 26451-2023"
$ws.Cells.Item(21, 7).Value = 2.2

# Row 22 becomes the "A 33852-2023" entry.
$ws.Cells.Item(22, 1).Value = "A 33852-2023"
$ws.Cells.Item(22, 2).Value = 45121
$ws.Cells.Item(22, 7).Value = 8.5

# Row 23 becomes the new "A 34013-2023" entry.
$ws.Cells.Item(23, 1).Value = "A 34013-2023"
$ws.Cells.Item(23, 2).Value = 45124
$ws.Cells.Item(23, 7).Value = 1.4

# Row 24 becomes "A 36814-2023" (previously row 23's data).
$ws.Cells.Item(24, 1).Value = "A 36814-2023"
$ws.Cells.Item(24, 2).Value = 45154
$ws.Cells.Item(24, 7).Value = 4.4
$ws.Rows.Item(24).RowHeight = 15

# Row 25 is the new row, holding the "A 36808-2023" entry (previously row 22's data).
$ws.Cells.Item(25, 1).Value = "A 36808-2023"
$ws.Cells.Item(25, 2).Value = 45154
$ws.Cells.Item(25, 3).Value = 45234
$ws.Cells.Item(25, 4).Value = "OKÄNT"
$ws.Cells.Item(25, 5).Value = "OKÄNT"
$ws.Cells.Item(25, 7).Value = 2.1
for ($c = 8; $c -le 17; $c++) {
    $ws.Cells.Item(25, $c).Value = 0
}
